# "Update framework tools." - the workbook was re-saved by a newer/different
# Excel build (Windows -> Mac path in the absPath hint, theme/version churn,
# etc.), but the only real content-level edits are:
#   1. The second sheet ("LString@built-in") was renamed to "LString@resource".
#   2. That sheet ("LString@resource") became the active/selected tab
#      (previously "LString" was the selected tab).

$wb = $excel.ActiveWorkbook

# 1) Rename the second worksheet.
$wsResource = $wb.Worksheets.Item(2)
$wsResource.Name = "LString@resource"

# 2) Make it the active sheet/tab.
$wsResource.Activate()
